# Actualización automática: 2026-01-13 10:18:13
# Updates the "Reporte TrendMicro" sheet: refreshes the host inventory rows
# (hosts renamed/renumbered, hostnames/IPs updated, version/release/tmxbc/
# ds_agent columns refreshed) and appends a new row 11 (MASTER02).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values (columns A..N) for data rows 2..11 after the update.
$rowData = @{
    2  = @("APP01",    "vm-prd-appps-bog-06",       "172.25.98.111",  "20.0.2", "29810.ol8", "20.0.2", "29810.ol8", "No instalado", "1.2.0.956",  "no activo", "active", "active", "active", "REVISAR")
    3  = @("APP02",    "vm-prd-appps-bog-07",       "172.25.96.193",  "20.0.2", "29810.ol8", "20.0.2", "29810.ol8", "No instalado", "1.2.0.956",  "no activo", "active", "active", "active", "REVISAR")
    4  = @("BATCH01",  "vm-prd-appps-bog-08",       "172.25.96.96",   "20.0.2", "29810.ol8", "20.0.2", "29810.ol8", "No instalado", "1.2.0.956",  "no activo", "active", "active", "active", "REVISAR")
    5  = @("BATCH02",  "vm-prd-appps-bog-09",       "172.25.99.10",   "20.0.2", "29810.ol8", "20.0.2", "29810.ol8", "No instalado", "1.2.0.956",  "no activo", "active", "active", "active", "REVISAR")
    6  = @("WS01",     "vm-prd-appps-bog-01-950287","172.25.97.44",   "20.0.2", "29760.ol8", "20.0.2", "29810.ol8", "1.2.0.1102",   "1.2.0.1102", "active",    "active", "active", "active", "REVISAR")
    7  = @("WS02",     "vm-prd-appps-bog-02",       "172.25.97.201",  "20.0.2", "29810.ol8", "20.0.2", "29810.ol8", "1.2.0.1253",   "1.2.0.1253", "active",    "active", "active", "active", "REVISAR")
    8  = @("BATCH03",  "vm-prd-appps-bog-10",       "172.25.99.176",  "20.0.2", "29810.ol8", "20.0.2", "29810.ol8", "No instalado", "1.2.0.956",  "no activo", "active", "active", "active", "REVISAR")
    9  = @("BATCH04",  "vm-prd-appps-bog-11",       "172.25.96.209",  "20.0.2", "29810.ol8", "20.0.2", "29810.ol8", "No instalado", "1.2.0.956",  "no activo", "active", "active", "active", "REVISAR")
    10 = @("MASTER01", "vm-prd-appps-bog-12",       "172.25.97.108",  "20.0.2", "29760.ol8", "20.0.2", "29760.ol8", "1.2.0.1102",   "1.2.0.1102", "active",    "active", "active", "active", "REVISAR")
    11 = @("MASTER02", "vm-prd-appps-bog-13",       "172.25.98.228",  "20.0.2", "29810.ol8", "20.0.2", "29810.ol8", "No instalado", "1.2.0.956",  "no activo", "active", "active", "active", "REVISAR")
}

# Row 11 is brand new: copy formatting (orange "REVISAR" fill in N10) down
# into N11 before writing its value.
$ws.Range("N10").Copy()
$ws.Range("N11").PasteSpecial(-4122)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

# Column widths: A 20->12, B 22->28 (ColumnWidth setter adds a ~0.8333
# padding offset vs. the stored OOXML <col width>, so subtract it back out).
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666
$ws.Columns.Item(2).ColumnWidth = 27.166666666666668
